# Updates the coin price/volume table on Sheet1 to match the latest
# coinranking.com snapshot (GitHub Actions scheduled refresh).
#
# Numeric-looking "Price" strings (single-dot values such as "521.33")
# are written through Formula with a leading apostrophe so Excel keeps
# them as text instead of silently parsing them into numbers - matching
# the workbook's existing inline-string "Price"/"Volume" columns. Values
# that already contain a thousands separator (e.g. "58.643.73") are never
# ambiguous, so a plain .Value assignment is enough for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.643.73"
$ws.Range("E2").Value = "  -0.63%  "

# Row 3
$ws.Range("D3").Value = "2.627.31"
$ws.Range("E3").Value = "  +0.20%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Formula = "'521.33"
$ws.Range("E5").Value = "  +1.67%  "

# Row 6
$ws.Range("D6").Formula = "'144.81"
$ws.Range("E6").Value = "  -2.43%  "

# Row 7
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
$ws.Range("D8").Formula = "'0.575"
$ws.Range("E8").Value = "  -0.24%  "

# Row 9
$ws.Range("D9").Value = "2.632.47"
$ws.Range("E9").Value = "  -0.51%  "

# Row 10
$ws.Range("D10").Formula = "'6.25"
$ws.Range("E10").Value = "  -3.40%  "

# Row 11
$ws.Range("E11").Value = "  -1.51%  "

# Row 12
$ws.Range("E12").Value = "  -1.78%  "

# Row 13
$ws.Range("E13").Value = "  -0.79%  "

# Row 14
$ws.Range("D14").Value = "3.086.30"
$ws.Range("E14").Value = "  +0.21%  "

# Row 15
$ws.Range("D15").Value = "58.683.20"
$ws.Range("E15").Value = "  -0.18%  "

# Row 16
$ws.Range("D16").Formula = "'20.74"
$ws.Range("E16").Value = "  -2.74%  "

# Row 17
$ws.Range("E17").Value = "  -2.35%  "

# Row 18
$ws.Range("D18").Value = "2.630.34"
$ws.Range("E18").Value = "  -0.16%  "

# Row 19
$ws.Range("D19").Formula = "'346.23"
$ws.Range("E19").Value = "  -0.26%  "

# Row 20
$ws.Range("E20").Value = "  -3.69%  "

# Row 21
$ws.Range("D21").Formula = "'10.19"
$ws.Range("E21").Value = "  -3.16%  "

# Row 22
$ws.Range("E22").Value = "  -1.28%  "

# Row 23
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("D24").Formula = "'61.62"
$ws.Range("E24").Value = "  +1.08%  "

# Row 25
$ws.Range("D25").Formula = "'0.414"
$ws.Range("E25").Value = "  -2.67%  "

# Row 26
$ws.Range("E26").Value = "  +1.00%  "

# Row 27
$ws.Range("E27").Value = "  +0.56%  "

# Row 28
$ws.Range("E28").Value = "  -3.96%  "

# Row 29
$ws.Range("E29").Value = "  -1.44%  "

# Row 30
$ws.Range("E30").Value = "  +0.16%  "

# Row 31
$ws.Range("E31").Value = "  -2.40%  "

# Row 32
$ws.Range("E32").Value = "  +0.49%  "

# Row 33
$ws.Range("D33").Formula = "'18.82"
$ws.Range("E33").Value = "  -1.62%  "

# Row 34
$ws.Range("D34").Formula = "'149.14"
$ws.Range("E34").Value = "  -0.24%  "

# Row 35
$ws.Range("D35").Formula = "'0.970"
$ws.Range("E35").Value = "  -2.35%  "

# Row 36
$ws.Range("D36").Formula = "'3.96"
$ws.Range("E36").Value = "  -2.36%  "

# Row 37
$ws.Range("D37").Formula = "'1.13"
$ws.Range("E37").Value = "  -1.38%  "

# Row 38
$ws.Range("D38").Formula = "'36.51"
$ws.Range("E38").Value = "  +0.51%  "

# Row 39
$ws.Range("E39").Value = "  -5.39%  "

# Row 40
$ws.Range("E40").Value = "  -1.71%  "

# Row 41
$ws.Range("D41").Formula = "'3.62"
$ws.Range("E41").Value = "  -1.16%  "

# Row 42
$ws.Range("D42").Formula = "'280.80"
$ws.Range("E42").Value = "  -3.91%  "

# Row 43
$ws.Range("E43").Value = "  +0.43%  "

# Row 44
$ws.Range("E44").Value = "  -1.64%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Formula = "'19.52"
$ws.Range("E45").Value = "  -1.44%  "

# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Formula = "'0.597"
$ws.Range("E46").Value = "  -4.22%  "

# Row 47
$ws.Range("D47").Formula = "'0.0521"
$ws.Range("E47").Value = "  -4.11%  "

# Row 48
$ws.Range("D48").Formula = "'10.30"
$ws.Range("E48").Value = "  +0.41%  "

# Row 49
$ws.Range("E49").Value = "  -1.79%  "

# Row 50
$ws.Range("D50").Value = "1.981.15"
$ws.Range("E50").Value = "  +0.59%  "

# Row 51
$ws.Range("D51").Formula = "'4.63"
$ws.Range("E51").Value = "  -3.37%  "
